$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(27.80092885869304, 23.67401530738431, 5.504233756908855, 15.96070750817679, 15.97227618470053, 21.21429875960476, 16.54237742217947),
    @(26.47546066207205, 22.71701502148927, 5.501159950652083, 8.125741681048671, 5.526617061661006, 10.74086696897962, 11.6700305017558),
    @(25.37732508547845, 21.73107982435348, 2.911879191335304, 5.520051518435407, 5.522992852459375, 13.34021356015267, 6.808995043443601),
    @(25.29361227171898, 21.72651166618357, 8.087175115391387, 0.3255479736593634, 2.922763035050882, 8.126862834164543, 6.808984613947648),
    @(24.5536477894986, 21.16739314693493, 2.911012591019551, 5.517024687163619, 5.519614635795364, 0.3253916112829664, 8.425483202850812),
    @(25.48003395504958, 21.87799954590978, 8.088029523928475, 15.93616035011776, 21.1699807137643, 8.129800689483552, 24.65785771812748),
    @(26.67111076250776, 22.92575579719268, 8.091985211070272, 10.7318283966353, 15.9576702852899, 8.135537216598438, 21.41286235947357),
    @(30.31496442985201, 25.89640409226083, 10.70116998032242, 8.14567959380604, 13.38020613335354, 18.62641935033201, 18.18176778406686),
    @(42.843675855665, 36.55080137981363, 18.60985231023347, 10.84467014114924, 16.13674213051041, 18.79369333917971, 29.71006928150475),
    @(53.87699545038097, 45.65839873882087, 8.190056381830818, 16.2436930055878, 16.26751482437832, 32.32737052271275, 23.24193205110446),
    @(72.5804014788592, 60.72209694555979, 10.91808667152936, 49.05275485384254, 54.6477345287984, 49.22986879958665, 34.97963250810057),
    @(86.40213106132725, 72.92364319467113, 29.76964580039761, 46.8174402247777, 35.90214440525371, 38.71785041179746, 72.14330730481157),
    @(87.71791673913042, 73.77269364477009, 32.48480096340825, 55.16116513169082, 35.94005890871207, 30.47901746057333, 58.65760577448611),
    @(86.79270311979064, 73.06747324374153, 19.03348079282802, 49.60277269723709, 55.28513598860674, 41.50726398690694, 68.76371182329824),
    @(86.84274244526338, 72.95957333816423, 27.0890568251119, 44.08621431488965, 52.5140321303595, 49.81966745320788, 63.69286401927418),
    @(87.14212753050239, 72.71054555915022, 32.47959572736801, 38.60071817175081, 80.45461010185365, 63.77581719121006, 55.2639084606892),
    @(73.68418821995799, 61.86341476461234, 29.56833288581934, 49.10538312192275, 74.00730745990953, 57.53460541159166, 73.46954391821666),
    @(50.55888141499607, 42.45255628575195, 42.44693680143635, 67.14611677595337, 83.59157517146177, 75.47946635022488, 101.444517027523),
    @(43.64130626330396, 36.93984296918191, 29.12150622417492, 58.72645299229748, 61.50524451300725, 50.80661482648864, 84.39737195904577),
    @(37.52231019108684, 31.64080988428539, 26.4198641619512, 74.57139146589144, 90.88918850993221, 61.30089626180726, 104.3650187312022),
    @(34.3272125581035, 29.29130979579332, 23.76237630525954, 61.00350471546056, 63.74941069217343, 31.88628082163538, 87.48115290037595),
    @(33.32943639711525, 28.24064199294787, 28.97870058576024, 42.37042223048101, 53.03566514087888, 39.79123260466888, 65.81121648880752),
    @(30.63281914849852, 26.02683563822602, 18.50675471111228, 23.85498114993836, 29.13410742278888, 29.1480136826575, 36.12965476153023),
    @(29.26345469573599, 25.03844374789078, 10.6963375859636, 15.97492962767688, 18.60433480806014, 10.76265863196949, 23.0576671353369)
)

$range = $ws.Range("B2:H25")
for ($i = 0; $i -lt $values.Count; $i++) {
    for ($j = 0; $j -lt $values[$i].Count; $j++) {
        $range.Cells.Item($i + 1, $j + 1).Value = $values[$i][$j]
    }
}
